$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Intents")

$ws.Range("B7").Value = "what's up, what’s happening, what’s new, whats up"
$ws.Range("B9").Value = "y, yes, okay, confirm, ok, yeah, sure "
$ws.Range("B10").Value = "n, no, nah, nope"
